$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously held a single header/data row (row 1) with the
# columns [reference, quantite, libelle]. The "nombre d'offres" column was
# invalid/missing, so the row is rebuilt starting at row 2 with a new
# "nombre d'offres" column inserted (col C), a similar second row, and a
# new wider "offre" detail row - then the whole 3-row block is duplicated
# once more below it.

# Clear out the old row 1 content entirely (it moves down / gets replaced).
$ws.Rows(1).ClearContents()

# Cells whose literal text looks like a number must be forced to Text so
# Excel stores them as shared strings ("3.0", "1", "2.0", "5", "8.0")
# instead of silently coercing them to numeric cells.
$textCells = @("B2","C2","B3","C3","B4","C4","B5","C5","B6","C6","B7","C7")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Cells.Item(2,1).Value = "ABCD7878"
$ws.Cells.Item(2,2).Value = "3.0"
$ws.Cells.Item(2,3).Value = "1"
$ws.Cells.Item(2,4).Value = "ABCDFGTRF89876*km"

# Row 3
$ws.Cells.Item(3,1).Value = "IOIOPO987"
$ws.Cells.Item(3,2).Value = "2.0"
$ws.Cells.Item(3,3).Value = "1"
$ws.Cells.Item(3,4).Value = "9876543*&*&ERTY"

# Row 4
$ws.Cells.Item(4,1).Value = "ABCHJUH"
$ws.Cells.Item(4,2).Value = "8.0"
$ws.Cells.Item(4,3).Value = "5"
$ws.Cells.Item(4,4).Value = "909ikokujyhtgt*"
$ws.Cells.Item(4,5).Value = "JKJKUHY/////\\\\\%^%gyvb"
$ws.Cells.Item(4,6).Value = "jkjkhjhh)))))"
$ws.Cells.Item(4,7).Value = "(((hjnmnmnmm####"
$ws.Cells.Item(4,8).Value = "hjhj.uiuiuisdksd"

# Row 5 (repeat of row 2's block)
$ws.Cells.Item(5,1).Value = "ABCD7878"
$ws.Cells.Item(5,2).Value = "3.0"
$ws.Cells.Item(5,3).Value = "1"
$ws.Cells.Item(5,4).Value = "ABCDFGTRF89876*km"

# Row 6 (repeat of row 3's block)
$ws.Cells.Item(6,1).Value = "IOIOPO987"
$ws.Cells.Item(6,2).Value = "2.0"
$ws.Cells.Item(6,3).Value = "1"
$ws.Cells.Item(6,4).Value = "9876543*&*&ERTY"

# Row 7 (repeat of row 4's block)
$ws.Cells.Item(7,1).Value = "ABCHJUH"
$ws.Cells.Item(7,2).Value = "8.0"
$ws.Cells.Item(7,3).Value = "5"
$ws.Cells.Item(7,4).Value = "909ikokujyhtgt*"
$ws.Cells.Item(7,5).Value = "JKJKUHY/////\\\\\%^%gyvb"
$ws.Cells.Item(7,6).Value = "jkjkhjhh)))))"
$ws.Cells.Item(7,7).Value = "(((hjnmnmnmm####"
$ws.Cells.Item(7,8).Value = "hjhj.uiuiuisdksd"
